$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stray spaces / punctuation from a handful of labels in column A
$ws.Range("A10").Value = "Evenwijdigelijnen"
$ws.Range("A12").Value = "Gelijkbenigerechthoekigedriehoek"
$ws.Range("A29").Value = "Opstaanderibben"
$ws.Range("A30").Value = "Opstaandezijvlakken"
$ws.Range("A43").Value = "Topruimtefiguur"
$ws.Range("A48").Value = "Vlakkefiguren"
$ws.Range("A48").Interior.ColorIndex = 6

# Fill in the previously-empty description cells in column B, matching the
# yellow header style used throughout column A where a B value is present
$ws.Range("B44").Value = "ruimtefiguur"
$ws.Range("A44").Interior.ColorIndex = 6

$ws.Range("B46").Value = "vlak figuur, zijden, hoekpunten, som, hoeken, graden."
$ws.Range("A46").Interior.ColorIndex = 6

$ws.Range("B49").Value = "Lijnen, vlak figuur, hoekpunten"
$ws.Range("A49").Interior.ColorIndex = 6

# B48 stays empty but picks up the same yellow highlight as its row's A cell
$ws.Range("B48").Interior.ColorIndex = 6

# Update the scrolled/selected cell shown when the workbook was saved
$ws.Range("A10").Select()
